$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Column C ("isMissing") currently holds the text string "FALSE" for every
# data row. Convert it to a real boolean FALSE value instead of text.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 301 }

$ws.Range("C2:C$lastRow").Value = $false
